$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Master Equipment Desc" for the W333 rows from "S2706 ASTUTE" to "Sonar"
$ws.Range("B3:B6").Value = "Sonar"

$ws.Range("B3:B6").Select()
